$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values: force text to avoid Excel auto-converting
# numeric-looking strings (e.g. "593.44") into actual numbers, which
# would change the stored cell type away from the original inline/shared string.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '61.113.84'
$ws.Range('E2').Value = '  +0.44%  '
Set-TextValue $ws.Range('D3') '2.928.43'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue $ws.Range('D5') '593.44'
$ws.Range('E5').Value = '  +1.13%  '
Set-TextValue $ws.Range('D6') '146.38'
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.02%  '
Set-TextValue $ws.Range('D9') '6.89'
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('E14').Value = '  -0.29%  '
Set-TextValue $ws.Range('D15') '3.412.58'
$ws.Range('E15').Value = '  +0.75%  '
Set-TextValue $ws.Range('D16') '61.084.63'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('E17').Value = '  -1.50%  '
Set-TextValue $ws.Range('D18') '2.923.44'
$ws.Range('E18').Value = '  +0.65%  '
Set-TextValue $ws.Range('D19') '432.20'
$ws.Range('E19').Value = '  +0.94%  '
Set-TextValue $ws.Range('D20') '13.46'
$ws.Range('E20').Value = '  -1.40%  '
Set-TextValue $ws.Range('D21') '0.684'
$ws.Range('E21').Value = '  +1.88%  '
Set-TextValue $ws.Range('D22') '7.09'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +5.73%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('E30').Value = '  +0.18%  '
Set-TextValue $ws.Range('D31') '7.09'
$ws.Range('E31').Value = '  -1.66%  '
Set-TextValue $ws.Range('D32') '26.51'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +1.34%  '
Set-TextValue $ws.Range('D34') '0.0₃0855'
$ws.Range('E34').Value = '  +2.28%  '
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('E36').Value = '  -0.47%  '
Set-TextValue $ws.Range('D37') '3.05'
$ws.Range('E37').Value = '  +3.44%  '
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('E39').Value = '  -1.70%  '
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  -1.69%  '
Set-TextValue $ws.Range('D42') '40.13'
$ws.Range('E42').Value = '  -4.97%  '
Set-TextValue $ws.Range('D43') '376.73'
$ws.Range('E43').Value = '  +1.31%  '
Set-TextValue $ws.Range('D44') '2.730.54'
$ws.Range('E44').Value = '  +2.71%  '
$ws.Range('E45').Value = '  +0.02%  '
Set-TextValue $ws.Range('D46') '130.26'
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('E47').Value = '  -0.07%  '
Set-TextValue $ws.Range('D48') '24.09'
$ws.Range('E48').Value = '  -3.49%  '
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('E50').Value = '  -2.95%  '
$ws.Range('E51').Value = '  +2.38%  '
